$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = 3
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = 0.003
$ws.Cells.Item(22, 5).Value = 'Regular'
$ws.Cells.Item(22, 6).Value = '<function relu at 0x11c4be9d8>'
$ws.Cells.Item(22, 7).Value = 0.9559000134468079
$ws.Cells.Item(22, 8).Value = 0.05790000036358833
$ws.Cells.Item(22, 9).Value = 0.04890000075101852
$ws.Cells.Item(22, 10).Value = 0.1738029420375824
$ws.Cells.Item(22, 11).Value = 6.665348529815674
$ws.Cells.Item(22, 12).Value = 0.05790000036358833
$ws.Cells.Item(22, 13).Value = 'logs/results_222.log'
$ws.Cells.Item(22, 14).Value = 'weights/model_222.ckpt'
$ws.Cells.Item(22, 15).Value = 'tb/222/non_robust'
$ws.Cells.Item(22, 16).Value = '(7.230026, 8.443772, 9.982576, 10.426454, 10.48337, 10.4655, 8.748777)'
$ws.Cells.Item(22, 17).Value = '(139.92761, 9.773078, 8.607171, 8.341426, 7.4237742, 7.5570107, 6.5512905, 8.353663)'

# Row 23
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = 3
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 0.003
$ws.Cells.Item(23, 5).Value = 'Regular'
$ws.Cells.Item(23, 6).Value = '<function relu at 0x11e21a9d8>'
$ws.Cells.Item(23, 7).Value = 0.9312000274658203
$ws.Cells.Item(23, 8).Value = 0.06909999996423721
$ws.Cells.Item(23, 9).Value = 0.1206000000238419
$ws.Cells.Item(23, 10).Value = 0.2596416771411896
$ws.Cells.Item(23, 11).Value = 7.456814289093018
$ws.Cells.Item(23, 12).Value = 0.06909999996423721
$ws.Cells.Item(23, 13).Value = 'logs/results_232.log'
$ws.Cells.Item(23, 14).Value = 'weights/model_232.ckpt'
$ws.Cells.Item(23, 15).Value = 'tb/232/non_robust'
$ws.Cells.Item(23, 16).Value = '(7.3609123, 7.8955326, 9.737668, 9.91078, 10.763038, 10.914308, 9.4478855)'
$ws.Cells.Item(23, 17).Value = '(140.87561, 8.819877, 8.838883, 8.998733, 7.921186, 7.691656, 7.321085, 9.176294)'

# Row 24
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = 3
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 0.003
$ws.Cells.Item(24, 5).Value = 'Regular'
$ws.Cells.Item(24, 6).Value = '<function relu at 0x111a129d8>'
$ws.Cells.Item(24, 7).Value = 0.9312999844551086
$ws.Cells.Item(24, 8).Value = 0.03830000013113022
$ws.Cells.Item(24, 9).Value = 0.009100000374019146
$ws.Cells.Item(24, 10).Value = 0.2439231872558594
$ws.Cells.Item(24, 11).Value = 7.289804458618164
$ws.Cells.Item(24, 12).Value = 0.03830000013113022
$ws.Cells.Item(24, 13).Value = 'logs/results_234.log'
$ws.Cells.Item(24, 14).Value = 'weights/model_234.ckpt'
$ws.Cells.Item(24, 15).Value = 'tb/234/non_robust'
$ws.Cells.Item(24, 16).Value = '(7.0172796, 7.6613173, 8.302304, 8.662146, 7.9462867, 9.3396435, 7.44935)'
$ws.Cells.Item(24, 17).Value = '(136.32106, 9.159889, 8.33736, 8.452992, 7.9244547, 7.2932935, 7.7801814, 9.626968)'
